# Updates cached market-price / profit figures on the Hyperion profit-tracking sheets.
# Values originate from a scheduled external data refresh (see commit message);
# this script just re-applies the refreshed numbers cell by cell.
$wb = $excel.ActiveWorkbook

$wsALC = $wb.Worksheets.Item("ALC")

# --- ALC sheet ---
$wsALC.Range("H33").Value = 4632.778
$wsALC.Range("I33").Value = 5314.1333
$wsALC.Range("K33").Value = 5314.1333
$wsALC.Range("M33").Value = -5085.1333
$wsALC.Range("H53").Value = 4592.76
$wsALC.Range("J53").Value = 10115.546
$wsALC.Range("L53").Value = 10115.546
$wsALC.Range("N53").Value = -11389.546
$wsALC.Range("H98").Value = 2617.4736
$wsALC.Range("I98").Value = 2207.389
$wsALC.Range("K98").Value = 2207.389
$wsALC.Range("M98").Value = -709.3890000000001
$wsALC.Range("H122").Value = 2617.4736
$wsALC.Range("I122").Value = 2207.389
$wsALC.Range("K122").Value = 6622.167
$wsALC.Range("M122").Value = -4172.167
$wsALC.Range("H125").Value = 2926.48
$wsALC.Range("J125").Value = 3283.65
$wsALC.Range("L125").Value = 29552.85
$wsALC.Range("N125").Value = -34472.85000000001
$wsALC.Range("H137").Value = 129451.71
$wsALC.Range("I137").Value = 356761.2
$wsALC.Range("J137").Value = 3168.6667
$wsALC.Range("K137").Value = 1070283.6
$wsALC.Range("L137").Value = 9506.000100000001
$wsALC.Range("M137").Value = -1067733.6
$wsALC.Range("N137").Value = -14606.0001
$wsALC.Range("H138").Value = 2668.673
$wsALC.Range("J138").Value = 4862.2383
$wsALC.Range("L138").Value = 14586.7149
$wsALC.Range("N138").Value = -24866.7149
$wsARM = $wb.Worksheets.Item("ARM")

# --- ARM sheet ---
$wsARM.Range("H2").Value = 1124.4117
$wsARM.Range("I2").Value = 775.25
$wsARM.Range("J2").Value = 1434.7778
$wsARM.Range("K2").Value = 775.25
$wsARM.Range("L2").Value = 1434.7778
$wsARM.Range("M2").Value = -662.25
$wsARM.Range("N2").Value = -1660.7778
$wsARM.Range("H97").Value = 10821.786
$wsARM.Range("I97").Value = 8643.521000000001
$wsARM.Range("K97").Value = 8643.521000000001
$wsARM.Range("M97").Value = -8147.521000000001
$wsARM.Range("H116").Value = 1124.4117
$wsARM.Range("I116").Value = 775.25
$wsARM.Range("J116").Value = 1434.7778
$wsARM.Range("K116").Value = 775.25
$wsARM.Range("L116").Value = 1434.7778
$wsARM.Range("M116").Value = 1518.75
$wsARM.Range("N116").Value = -6022.7778
$wsARM.Range("H122").Value = 2969.074
$wsARM.Range("I122").Value = 2753.15
$wsARM.Range("K122").Value = 8259.450000000001
$wsARM.Range("M122").Value = -5809.450000000001
$wsARM.Range("H132").Value = 2916.2666
$wsARM.Range("I132").Value = 2694
$wsARM.Range("K132").Value = 8082
$wsARM.Range("M132").Value = -5552
$wsBSM = $wb.Worksheets.Item("BSM")

# --- BSM sheet ---
$wsBSM.Range("H3").Value = 1124.4117
$wsBSM.Range("I3").Value = 775.25
$wsBSM.Range("J3").Value = 1434.7778
$wsBSM.Range("K3").Value = 775.25
$wsBSM.Range("L3").Value = 1434.7778
$wsBSM.Range("M3").Value = -661.25
$wsBSM.Range("N3").Value = -1662.7778
$wsBSM.Range("H99").Value = 3042.0952
$wsBSM.Range("I99").Value = 1891.5834
$wsBSM.Range("J99").Value = 4576.1113
$wsBSM.Range("K99").Value = 1891.5834
$wsBSM.Range("L99").Value = 4576.1113
$wsBSM.Range("M99").Value = -393.5834
$wsBSM.Range("N99").Value = -7572.1113
$wsBSM.Range("H105").Value = 1668.75
$wsBSM.Range("I105").Value = 1668.75
$wsBSM.Range("K105").Value = 1668.75
$wsBSM.Range("M105").Value = 78.25
$wsBSM.Range("H134").Value = 3202.5715
$wsBSM.Range("I134").Value = 1627.9688
$wsBSM.Range("K134").Value = 4883.9064
$wsBSM.Range("M134").Value = -2348.9064
$wsCRP = $wb.Worksheets.Item("CRP")

# --- CRP sheet ---
$wsCRP.Range("H31").Value = 4213.52
$wsCRP.Range("I31").Value = 1346.5883
$wsCRP.Range("J31").Value = 4800.7227
$wsCRP.Range("K31").Value = 1346.5883
$wsCRP.Range("L31").Value = 4800.7227
$wsCRP.Range("M31").Value = -1051.5883
$wsCRP.Range("N31").Value = -5390.7227
$wsCRP.Range("H34").Value = 4213.52
$wsCRP.Range("I34").Value = 1346.5883
$wsCRP.Range("J34").Value = 4800.7227
$wsCRP.Range("K34").Value = 1346.5883
$wsCRP.Range("L34").Value = 4800.7227
$wsCRP.Range("M34").Value = -1144.5883
$wsCRP.Range("N34").Value = -5204.7227
$wsCRP.Range("H86").Value = 10994.066
$wsCRP.Range("J86").Value = 12141.857
$wsCRP.Range("L86").Value = 12141.857
$wsCRP.Range("N86").Value = -14387.857
$wsCRP.Range("H89").Value = 10994.066
$wsCRP.Range("J89").Value = 12141.857
$wsCRP.Range("L89").Value = 60709.285
$wsCRP.Range("N89").Value = -71941.285
$wsCRP.Range("H94").Value = 1235.2727
$wsCRP.Range("I94").Value = 869
$wsCRP.Range("J94").Value = 1271.9
$wsCRP.Range("K94").Value = 869
$wsCRP.Range("L94").Value = 1271.9
$wsCRP.Range("M94").Value = -418
$wsCRP.Range("N94").Value = -2173.9
$wsCRP.Range("H99").Value = 3196.2727
$wsCRP.Range("I99").Value = 2451.7144
$wsCRP.Range("K99").Value = 2451.7144
$wsCRP.Range("M99").Value = -953.7143999999998
$wsCRP.Range("H105").Value = 2224.7144
$wsCRP.Range("I105").Value = 1928.8334
$wsCRP.Range("K105").Value = 1928.8334
$wsCRP.Range("M105").Value = -181.8334
$wsCRP.Range("H107").Value = 45456310
$wsCRP.Range("I107").Value = 1693.3529
$wsCRP.Range("J107").Value = 200002020
$wsCRP.Range("K107").Value = 1693.3529
$wsCRP.Range("L107").Value = 200002020
$wsCRP.Range("M107").Value = 226.6470999999999
$wsCRP.Range("N107").Value = -200005860
$wsCRP.Range("H126").Value = 3196.2727
$wsCRP.Range("I126").Value = 2451.7144
$wsCRP.Range("K126").Value = 7355.1432
$wsCRP.Range("M126").Value = -4885.1432
$wsCRP.Range("H132").Value = 30137.312
$wsCRP.Range("I132").Value = 2453.5334
$wsCRP.Range("K132").Value = 7360.600199999999
$wsCRP.Range("M132").Value = -4830.600199999999
$wsCRP.Range("H134").Value = 3519.5862
$wsCRP.Range("I134").Value = 3405.8
$wsCRP.Range("K134").Value = 10217.4
$wsCRP.Range("M134").Value = -7682.400000000001
$wsCUL = $wb.Worksheets.Item("CUL")

# --- CUL sheet ---
$wsCUL.Range("H37").Value = 46454.89
$wsCUL.Range("J37").Value = 46454.89
$wsCUL.Range("L37").Value = 139364.67
$wsCUL.Range("N37").Value = -139588.67
$wsCUL.Range("H56").Value = 50004920
$wsCUL.Range("I56").Value = 50004920
$wsCUL.Range("K56").Value = 50004920
$wsCUL.Range("M56").Value = -50004390
$wsCUL.Range("H97").Value = 2622.3845
$wsCUL.Range("J97").Value = 497
$wsCUL.Range("L97").Value = 1491
$wsCUL.Range("N97").Value = -2483
$wsCUL.Range("H117").Value = 23813056
$wsCUL.Range("I117").Value = 55560070
$wsCUL.Range("J117").Value = 2796.875
$wsCUL.Range("K117").Value = 166680210
$wsCUL.Range("L117").Value = 8390.625
$wsCUL.Range("M117").Value = -166676768
$wsCUL.Range("N117").Value = -15274.625
$wsCUL.Range("H140").Value = 1936
$wsCUL.Range("I140").Value = 1837.7142
$wsCUL.Range("K140").Value = 5513.142599999999
$wsCUL.Range("M140").Value = -333.1425999999992
$wsGSM = $wb.Worksheets.Item("GSM")

# --- GSM sheet ---
$wsGSM.Range("H15").Value = 15662.333
$wsGSM.Range("I15").Value = 12000
$wsGSM.Range("J15").Value = 17493.5
$wsGSM.Range("K15").Value = 12000
$wsGSM.Range("L15").Value = 17493.5
$wsGSM.Range("M15").Value = -11712
$wsGSM.Range("N15").Value = -18069.5
$wsGSM.Range("H81").Value = 15662.333
$wsGSM.Range("I81").Value = 12000
$wsGSM.Range("J81").Value = 17493.5
$wsGSM.Range("K81").Value = 12000
$wsGSM.Range("L81").Value = 17493.5
$wsGSM.Range("M81").Value = -11002
$wsGSM.Range("N81").Value = -19489.5
$wsGSM.Range("H84").Value = 15662.333
$wsGSM.Range("I84").Value = 12000
$wsGSM.Range("J84").Value = 17493.5
$wsGSM.Range("K84").Value = 36000
$wsGSM.Range("L84").Value = 52480.5
$wsGSM.Range("M84").Value = -31008
$wsGSM.Range("N84").Value = -62464.5
$wsGSM.Range("H122").Value = 745480
$wsGSM.Range("I122").Value = 893076
$wsGSM.Range("J122").Value = 7500
$wsGSM.Range("K122").Value = 2679228
$wsGSM.Range("L122").Value = 22500
$wsGSM.Range("M122").Value = -2676778
$wsGSM.Range("N122").Value = -27400
$wsGSM.Range("H132").Value = 4294.846
$wsGSM.Range("I132").Value = 3398
$wsGSM.Range("K132").Value = 10194
$wsGSM.Range("M132").Value = -7664
$wsLTW = $wb.Worksheets.Item("LTW")

# --- LTW sheet ---
$wsLTW.Range("H22").Value = 129350.57
$wsLTW.Range("J22").Value = 3742.3333
$wsLTW.Range("L22").Value = 3742.3333
$wsLTW.Range("N22").Value = -4332.3333
$wsLTW.Range("H27").Value = 129350.57
$wsLTW.Range("J27").Value = 3742.3333
$wsLTW.Range("L27").Value = 3742.3333
$wsLTW.Range("N27").Value = -3956.3333
$wsLTW.Range("H40").Value = 7158.6
$wsLTW.Range("I40").Value = 5938.1
$wsLTW.Range("K40").Value = 5938.1
$wsLTW.Range("M40").Value = -5802.1
$wsLTW.Range("H68").Value = 3106.5715
$wsLTW.Range("I68").Value = 2486
$wsLTW.Range("J68").Value = 3934
$wsLTW.Range("K68").Value = 2486
$wsLTW.Range("L68").Value = 3934
$wsLTW.Range("M68").Value = -1737
$wsLTW.Range("N68").Value = -5432
$wsLTW.Range("H71").Value = 3106.5715
$wsLTW.Range("I71").Value = 2486
$wsLTW.Range("J71").Value = 3934
$wsLTW.Range("K71").Value = 12430
$wsLTW.Range("L71").Value = 19670
$wsLTW.Range("M71").Value = -8686
$wsLTW.Range("N71").Value = -27158
$wsLTW.Range("H136").Value = 37963.863
$wsLTW.Range("I136").Value = 55108.105
$wsLTW.Range("K136").Value = 165324.315
$wsLTW.Range("M136").Value = -162774.315
$wsWVR = $wb.Worksheets.Item("WVR")

# --- WVR sheet ---
$wsWVR.Range("H62").Value = 7165.074
$wsWVR.Range("J62").Value = 9559.843999999999
$wsWVR.Range("L62").Value = 9559.843999999999
$wsWVR.Range("N62").Value = -10807.844
$wsWVR.Range("H65").Value = 7165.074
$wsWVR.Range("J65").Value = 9559.843999999999
$wsWVR.Range("L65").Value = 47799.21999999999
$wsWVR.Range("N65").Value = -54039.21999999999
$wsWVR.Range("H107").Value = 5254.2666
$wsWVR.Range("I107").Value = 3736.4614
$wsWVR.Range("J107").Value = 15120
$wsWVR.Range("K107").Value = 11209.3842
$wsWVR.Range("L107").Value = 45360
$wsWVR.Range("M107").Value = -9289.3842
$wsWVR.Range("N107").Value = -49200
$wsWVR.Range("H109").Value = 69987.664
$wsWVR.Range("J109").Value = 69987.664
$wsWVR.Range("L109").Value = 69987.664
$wsWVR.Range("N109").Value = -72761.664
$wsWVR.Range("H132").Value = 33212.656
$wsWVR.Range("I132").Value = 8824.1875
$wsWVR.Range("J132").Value = 57601.125
$wsWVR.Range("K132").Value = 26472.5625
$wsWVR.Range("L132").Value = 172803.375
$wsWVR.Range("M132").Value = -23942.5625
$wsWVR.Range("N132").Value = -177863.375
$wsWVR.Range("H136").Value = 2619.2856
$wsWVR.Range("I136").Value = 1900.375
$wsWVR.Range("J136").Value = 4919.8
$wsWVR.Range("K136").Value = 5701.125
$wsWVR.Range("L136").Value = 14759.4
$wsWVR.Range("M136").Value = -3151.125
$wsWVR.Range("N136").Value = -19859.4
